# Update Name of Algo
# Applies updated KNN-imputed values to the result_data_KNN worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.235000000000001
$ws.Range("C3").Value = -11.981
$ws.Range("A4").Value = -21.574
$ws.Range("B4").Value = 6.525
$ws.Range("C4").Value = -12.081
$ws.Range("E4").Value = 13.262
$ws.Range("B5").Value = 6.379999999999999
$ws.Range("D5").Value = -8.328999999999999
$ws.Range("A6").Value = -21.259
$ws.Range("B6").Value = 6.822
$ws.Range("A7").Value = -21.179
$ws.Range("A8").Value = -21.379
$ws.Range("B8").Value = 6.361000000000001
$ws.Range("C9").Value = -11.775
$ws.Range("C11").Value = -12.642
$ws.Range("E12").Value = 12.959
$ws.Range("C14").Value = -11.784
$ws.Range("A16").Value = -20.95
$ws.Range("B16").Value = 6.622
$ws.Range("E16").Value = 13.408
$ws.Range("E17").Value = 12.786
$ws.Range("C18").Value = -12.634
$ws.Range("A20").Value = -22.185
$ws.Range("D20").Value = -8.280000000000001
$ws.Range("E20").Value = 13.154
$ws.Range("A21").Value = -21.14
$ws.Range("B22").Value = 6.386000000000001
$ws.Range("C25").Value = -12.43
$ws.Range("E25").Value = 13.078
